$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert ".chr1" into the H-column formulas (sj.bed.gz file name) for rows 2-5.
# H2 holds its own (non-shared) formula; H3:H5 share one formula (relative to H3).
$ws.Range("H2").Formula = '=CONCATENATE(B2,".chr1.sj.bed.gz")'
$ws.Range("H3:H5").Formula = '=CONCATENATE(B3,".chr1.sj.bed.gz")'

# Update the current selection on the sheet to H3:H5 (active cell H3).
$ws.Range("H3:H5").Select()

# Reposition the workbook window (best effort - not all headless hosts persist this).
$win = $excel.ActiveWindow
$win.Left = 4460
$win.Top = 2820
